# Minuta de Reunion template: tighten up heading spacing and page margins.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Heading 1") {
        # w:spacing w:before="80" -> w:before="120"  (80 twips = 4pt, 120 twips = 6pt)
        $p.SpaceBefore = 6
    }
    elseif ($styleName -eq "Heading 2") {
        # w:spacing w:before="80" -> w:before="100"  (80 twips = 4pt, 100 twips = 5pt)
        $p.SpaceBefore = 5
    }
}

# Page margins / header-footer distance (values are twips/20 = points)
$sec = $d.Sections.First
$ps = $sec.PageSetup
$ps.TopMargin = 76.55      # 1418 -> 1531 twips
$ps.BottomMargin = 76.55   # 1418 -> 1531 twips
$ps.HeaderDistance = 39.7  # 737 -> 794 twips
$ps.FooterDistance = 34.0  # 624 -> 680 twips
